#
# Distinguish between roles Student and GradStudent, because Student cannot
# generate a prototype and a GradStudent can.
#
# Adds a new "GradStudent" role value and applies it to the accounts that
# should get the extra (graduate) privilege: Stef, Lloyd and Debbie. Lloyd's
# account previously had no entry in the third "Role" column (I), so that
# cell is newly populated; Stef's and Debbie's existing "Student" entries in
# column I are upgraded to "GradStudent".
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Acc_Stef (row 14), Acc_Lloyd (row 15), Acc_Debbie (row 18)
$ws.Range("I14").Value = "GradStudent"
$ws.Range("I15").Value = "GradStudent"
$ws.Range("I18").Value = "GradStudent"

# Reflect the selection left active on the sheet after the edit.
[void]$ws.Range("G20").Select()
